$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 160-1827-1-ND / D1-12 / 12  ->  1568-1800-ND / D1-6 / 6
$ws.Range("A2").Value = "1568-1800-ND"
$ws.Range("B2").Value = "D1-6"
$ws.Range("C2").Value = 6

# Row 3: RR08P100DCT-ND / R1-12 / 12  ->  RR08P15.0KDCT-ND / R7-12 / 6
$ws.Range("A3").Value = "RR08P15.0KDCT-ND"
$ws.Range("B3").Value = "R7-12"
$ws.Range("C3").Value = 6

# Row 4: RNCP0805FTD20K0CT-ND / R13-18 / 6  ->  RNCP0805FTD20K0CT-ND / R1-6 / 6
$ws.Range("A4").Value = "RNCP0805FTD20K0CT-ND"
$ws.Range("B4").Value = "R1-6"
$ws.Range("C4").Value = 6

# Row 5: new part row - 732-5309-ND / P1 / 1, bold Arial 9pt label
$ws.Range("A5").Value = "732-5309-ND"
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").Font.Size = 9
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("B5").Value = "P1"
$ws.Range("C5").Value = 1

# Update the selected cell to match the saved workbook's cursor position
$ws.Range("F8").Select()
